# ShadeSwapCalculation.xlsx -- add "Staking" worksheet with staking-reward
# calculations (commit: "add staking calculation into excel.")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the two existing sheets' selections (cosmetic cursor moves
#    that were part of the author's saved session).
# ---------------------------------------------------------------------
$wsAmount = $wb.Worksheets.Item(1)
[void]$wsAmount.Range("I13").Select()

$wsLiquidity = $wb.Worksheets.Item(2)
[void]$wsLiquidity.Range("F17").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Staking" sheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Staking"

# ---------------------------------------------------------------------
# 3. Column widths (approximated to this engine's nearest achievable
#    "characters" granularity).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.142857142857142
$ws.Columns.Item(2).ColumnWidth = 27.714285714285715
$ws.Columns.Item(3).ColumnWidth = 27.714285714285715
$ws.Columns.Item(4).ColumnWidth = 34.714285714285715
$ws.Columns.Item(5).ColumnWidth = 20.714285714285715
$ws.Columns.Item(6).ColumnWidth = 30.714285714285715
$ws.Columns.Item(7).ColumnWidth = 20.714285714285715
$ws.Columns.Item(8).ColumnWidth = 20.714285714285715
$ws.Columns.Item(9).ColumnWidth = 22.142857142857142

# ---------------------------------------------------------------------
# 4. Header row + row labels -- written in the exact order the original
#    author entered them so new shared-string indices line up with the
#    source workbook (B1, G1, A1, A2, A4, A3, A5, D1, E1, I1, C1, A6,
#    F1, H1).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "LP Token"
$ws.Range("G1").Value = "Denominator"
$ws.Range("A1").Value = "Stakers"
$ws.Range("A2").Value = "Staker A"
$ws.Range("A4").Value = "Staker C"
$ws.Range("A3").Value = "Staker B"
$ws.Range("A5").Value = "Total LP"
$ws.Range("D1").Value = "Last Claimable Block Height"
$ws.Range("E1").Value = "NextClaim Time "
$ws.Range("I1").Value = "Claimable Reward"
$ws.Range("C1").Value = "% Staking "
$ws.Range("A6").Value = "Reward Amount"
$ws.Range("F1").Value = "Last Time - Current Time in Mils"
$ws.Range("H1").Value = "Total Available Reward"

# ---------------------------------------------------------------------
# 5. Raw numeric inputs.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 100000
$ws.Range("B3").Value = 2000000
$ws.Range("B4").Value = 4560000000
$ws.Range("B6").Value = 3450000000000

$ws.Range("D2").Value = 1656480000
$ws.Range("E2").Value = 1656480524
$ws.Range("D3").Value = 1656480000
$ws.Range("E3").Value = 1656480524
$ws.Range("D4").Value = 1656480000
$ws.Range("E4").Value = 1656480524

# ---------------------------------------------------------------------
# 6. Formulas (ranges written as one Formula assignment become proper
#    OOXML shared formulas, matching the source file).
# ---------------------------------------------------------------------
$ws.Range("B5").Formula = "=SUM(B2:B4)"

$ws.Range("C2").Formula = "=B2/`$B5"
$ws.Range("C3").Formula = "=B3/`$B5"
$ws.Range("C4").Formula = "=B4/`$B5"

$ws.Range("F2").Formula = "=E2-D2"
$ws.Range("F3:F4").Formula = "=E3-D3"

$ws.Range("G2").Formula = "=24*60*60*1000"
$ws.Range("G3:G4").Formula = "=24*60*60*1000"

$ws.Range("H2").Formula = "=(B6/G2) *F2"
$ws.Range("H3").Formula = "=(B6/G3) *F3"
$ws.Range("H4").Formula = "=(B6/G4) *F4"

$ws.Range("I2").Formula = "=INT(H2*C2)"
$ws.Range("I3:I4").Formula = "=INT(H3*C3)"

# ---------------------------------------------------------------------
# 7. Styling: thick-bordered "Check Cell" separator row beneath the
#    totals, and the bottom "Total Available Reward" row.
# ---------------------------------------------------------------------
$ws.Range("C5:H5").Style = "Check Cell"
$ws.Range("C6:H6").Style = "Check Cell"

# ---------------------------------------------------------------------
# 8. Row heights for the thick-bordered rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(7).RowHeight = 15.75

# ---------------------------------------------------------------------
# 9. Selection / active cell on the new sheet.
# ---------------------------------------------------------------------
[void]$ws.Range("F18").Select()
